$wb = $excel.ActiveWorkbook

# Insert the new "Review Data" worksheet right after "README" (before "Procedure").
$readme = $wb.Worksheets.Item("README")
$reviewData = $wb.Worksheets.Add($null, $readme)
$reviewData.Name = "Review Data"

# Populate the cells in the same order the values were originally typed in so
# that shared-string ids line up with the authored workbook.
$reviewData.Range("A1").Value = "Spec"
$reviewData.Range("B1").Value = "Test"
$reviewData.Range("C1").Value = "Comments"
$reviewData.Range("A8").Value = "Manual:"
$reviewData.Range("A9").Value = "header property shown as window title"
$reviewData.Range("A10").Value = "button property results in `"top`" button on ctrl panel"
$reviewData.Range("A11").Value = "window property can name alternate window"
$reviewData.Range("A20").Value = "omit posting method??"
$reviewData.Range("A13").Value = "button response saved in recipe data key"
$reviewData.Range("A14").Value = "gui position/scale control dialog location/size"
$reviewData.Range("A16").Value = "configuration stored @ configuration key"
$reviewData.Range("B9").Value = "x"
$reviewData.Range("A21").Value = "?? Don't see prompt show up in screen shot"
$reviewData.Range("A15").Value = "units conversion between recipe data units and configured `"view`" units"
$reviewData.Range("C16").Value = "in a string valued recipe data editable with s88Set"
$reviewData.Range("A12").Value = "posting method populates window"
$reviewData.Range("B12").Value = "specify default, takes same params"
$reviewData.Range("A22").Value = "?? Omit display mode, make secondary data a tab in same window"
$reviewData.Range("A17").Value = "values read from recipe data (value with units)"

# Column widths (best match achievable for 49 / 28.7109375 / 43.140625 characters).
$reviewData.Columns.Item(1).ColumnWidth = 48.166666666666664
$reviewData.Columns.Item(2).ColumnWidth = 27.833333333333332
$reviewData.Columns.Item(3).ColumnWidth = 42.333333333333336

# Make "Review Data" the active/visible tab with the same selection/scroll state
# as the authored workbook.
$reviewData.Activate() | Out-Null
$reviewData.Range("A18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
